$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated score values (column D) ---
$ws.Range("D4").Value = 5
$ws.Range("D13").Value = 5
$ws.Range("D14").Value = 4.857758620689655
$ws.Range("D41").Value = 4.584905660377358
$ws.Range("D47").Value = 4.869565217391305

# --- New rows 52-55: name + zero metrics + D=1 ---
$ws.Range("A52").Value = "Павликов Илья Сергеевич "
$ws.Range("B52").Value = 0
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 1

$ws.Range("A53").Value = "Верле Каролина Валерьевна (Обучение 2)"
$ws.Range("B53").Value = 0
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 1

$ws.Range("A54").Value = "Довыдович Алиса Станиславовна "
$ws.Range("B54").Value = 0
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 1

$ws.Range("A55").Value = " Шептунова Софья Денисовна"
$ws.Range("B55").Value = 0
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 1
